$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the old row 8 (which held "extr1"), pushing
# the old "extr1".."extr8" rows down from rows 8-15 to rows 10-17 while
# keeping their existing shared-string references for column B intact.
$ws.Rows("8:9").Insert()

# Copy the direct cell formatting (bold font, border, centered alignment)
# used on the indexed column (A) of the existing data rows onto the two
# newly inserted rows.
$ws.Range("A10").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)

# New row 8 -> name "line7"
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $false

# New row 9 -> name "line8"
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Renumber / update the shifted "extr" rows (now at rows 10-17), leaving
# their name (column B) shared-string references untouched.
$ws.Range("A10").Value = 8
$ws.Range("E10").Value = $true

$ws.Range("A11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("A12").Value = 10

$ws.Range("A13").Value = 11

$ws.Range("A14").Value = 12

$ws.Range("A15").Value = 13

$ws.Range("A16").Value = 14

$ws.Range("A17").Value = 15
